# Update cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.597.60"
$ws.Range("E2").Value = "'  +1.46%  "
$ws.Range("D3").Value = "'1.600.94"
$ws.Range("E3").Value = "'  +1.48%  "
$ws.Range("E4").Value = "'  +0.49%  "
$ws.Range("E5").Value = "'  +0.19%  "
$ws.Range("E6").Value = "'  +0.26%  "
$ws.Range("E7").Value = "'  +0.52%  "
$ws.Range("D8").Value = "'26.74"
$ws.Range("E8").Value = "'  +2.60%  "
$ws.Range("D9").Value = "'0.250"
$ws.Range("E9").Value = "'  +1.35%  "
$ws.Range("E10").Value = "'  +1.26%  "
$ws.Range("E11").Value = "'  +0.98%  "
$ws.Range("D12").Value = "'1.829.66"
$ws.Range("E12").Value = "'  +1.47%  "
$ws.Range("D13").Value = "'1.600.98"
$ws.Range("E13").Value = "'  +1.51%  "
$ws.Range("D14").Value = "'29.602.03"
$ws.Range("E14").Value = "'  +1.38%  "
$ws.Range("E15").Value = "'  +3.03%  "
$ws.Range("E16").Value = "'  +0.94%  "
$ws.Range("D18").Value = "'241.51"
$ws.Range("E18").Value = "'  +1.40%  "
$ws.Range("D19").Value = "'7.64"
$ws.Range("E19").Value = "'  +2.58%  "
$ws.Range("D20").Value = "'0.0₃0693"
$ws.Range("E21").Value = "'  +0.48%  "
$ws.Range("E22").Value = "'  -0.09%  "
$ws.Range("E23").Value = "'  +0.51%  "
$ws.Range("E24").Value = "'  -0.89%  "
$ws.Range("D25").Value = "'155.14"
$ws.Range("E25").Value = "'  +1.20%  "
$ws.Range("D26").Value = "'15.34"
$ws.Range("E27").Value = "'  +0.42%  "
$ws.Range("E28").Value = "'  +1.27%  "
$ws.Range("E29").Value = "'  +0.49%  "
$ws.Range("E30").Value = "'  +2.66%  "
$ws.Range("E31").Value = "'  +0.33%  "
$ws.Range("E32").Value = "'  +0.43%  "
$ws.Range("E33").Value = "'  +2.61%  "
$ws.Range("D34").Value = "'1.423.64"
$ws.Range("E34").Value = "'  -0.08%  "
$ws.Range("B35").Value = "'LidoDAOToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = "'  +2.40%  "
$ws.Range("B36").Value = "'MXToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'2.88"
$ws.Range("E36").Value = "'  +5.42%  "
$ws.Range("D37").Value = "'1.02"
$ws.Range("E37").Value = "'  -2.17%  "
$ws.Range("D38").Value = "'2.30"
$ws.Range("E38").Value = "'  +0.45%  "
$ws.Range("E39").Value = "'  +2.51%  "
$ws.Range("D40").Value = "'0.544"
$ws.Range("E40").Value = "'  +2.98%  "
$ws.Range("E41").Value = "'  +0.93%  "
$ws.Range("E42").Value = "'  +5.18%  "
$ws.Range("D43").Value = "'54.24"
$ws.Range("E43").Value = "'  +2.37%  "
$ws.Range("E44").Value = "'  +2.57%  "
$ws.Range("E45").Value = "'  +0.51%  "
$ws.Range("D46").Value = "'0.988"
$ws.Range("E46").Value = "'  +16.41%  "
$ws.Range("D47").Value = "'66.38"
$ws.Range("E47").Value = "'  +3.10%  "
$ws.Range("E48").Value = "'  -0.62%  "
$ws.Range("D49").Value = "'1.740.33"
$ws.Range("E49").Value = "'  +1.44%  "
$ws.Range("D50").Value = "'85.96"
$ws.Range("E50").Value = "'  +0.30%  "
$ws.Range("E51").Value = "'  +2.94%  "
